$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 985
$ws.Range("F6").Value = 3223
$ws.Range("F9").Value = 1144
$ws.Range("F14").Value = 75
$ws.Range("F15").Value = 634
$ws.Range("F16").Value = 1417
$ws.Range("F17").Value = 1417
$ws.Range("F18").Value = 248
$ws.Range("F21").Value = 281
$ws.Range("F23").Value = 478
$ws.Range("F24").Value = 25547
$ws.Range("F25").Value = 25548
$ws.Range("F28").Value = 16164
$ws.Range("F29").Value = 16164
$ws.Range("F30").Value = 374
$ws.Range("F33").Value = 897
$ws.Range("F34").Value = 187
$ws.Range("F36").Value = 441
$ws.Range("F37").Value = 1120
$ws.Range("F38").Value = 5244
$ws.Range("F39").Value = 650
$ws.Range("F40").Value = 385
$ws.Range("F42").Value = 294
$ws.Range("F45").Value = 34

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F13").Value = 1786
$ws.Range("F17").Value = 383
$ws.Range("F34").Value = 848
$ws.Range("F35").Value = 477
$ws.Range("F36").Value = 3
$ws.Range("F37").Value = 59

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 500
$ws.Range("F6").Value = 492

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 500
$ws.Range("F6").Value = 985
$ws.Range("F14").Value = 492
$ws.Range("F17").Value = 1786
$ws.Range("F21").Value = 75
$ws.Range("F22").Value = 634
$ws.Range("F23").Value = 1417
$ws.Range("F24").Value = 1417
$ws.Range("F25").Value = 248
$ws.Range("F29").Value = 383
$ws.Range("F30").Value = 281
$ws.Range("F31").Value = 478
$ws.Range("F33").Value = 25548
$ws.Range("F36").Value = 16165
$ws.Range("F37").Value = 374
$ws.Range("F38").Value = 897
$ws.Range("F40").Value = 187
$ws.Range("F43").Value = 441
$ws.Range("F44").Value = 5244
$ws.Range("F46").Value = 650
$ws.Range("F47").Value = 477
$ws.Range("F48").Value = 59
$ws.Range("F49").Value = 59
$ws.Range("F50").Value = 294
$ws.Range("F54").Value = 34
